$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update scattered odds values in existing rows ---
$ws.Range("U2").Value2 = 1.62
$ws.Range("U3").Value2 = 1.75
$ws.Range("V4").Value2 = 1.7
$ws.Range("V5").Value2 = 1.67
$ws.Range("BD5").Value2 = 151
$ws.Range("N10").Value2 = 10
$ws.Range("Q10").Value2 = 2
$ws.Range("R10").Value2 = 1.85
$ws.Range("N18").Value2 = 15
$ws.Range("Q18").Value2 = 1.62
$ws.Range("R18").Value2 = 2.25
$ws.Range("Q19").Value2 = 1.7
$ws.Range("R19").Value2 = 2.1
$ws.Range("V21").Value2 = 1.63
$ws.Range("V22").Value2 = 1.72
$ws.Range("U23").Value2 = 1.92
$ws.Range("V23").Value2 = 1.77
$ws.Range("V24").Value2 = 1.54
$ws.Range("U25").Value2 = 2.37
$ws.Range("V25").Value2 = 1.5
$ws.Range("U33").Value2 = 1.83
$ws.Range("V33").Value2 = 1.83
$ws.Range("U34").Value2 = 1.95
$ws.Range("V34").Value2 = 1.8
$ws.Range("O46").Value2 = 1.29
$ws.Range("P46").Value2 = 3.5
$ws.Range("Q46").Value2 = 2
$ws.Range("R46").Value2 = 1.85
$ws.Range("M48").Value2 = 1.03
$ws.Range("O48").Value2 = 1.14
$ws.Range("M49").Value2 = 1.04
$ws.Range("O49").Value2 = 1.25
$ws.Range("M50").Value2 = 1.03
$ws.Range("O50").Value2 = 1.22
$ws.Range("M51").Value2 = 1.03
$ws.Range("O51").Value2 = 1.25
$ws.Range("U52").Value2 = 1.92
$ws.Range("V52").Value2 = 1.77
$ws.Range("U53").Value2 = 1.87
$ws.Range("V53").Value2 = 1.87
$ws.Range("M54").Value2 = 1.07
$ws.Range("O54").Value2 = 1.33
$ws.Range("U54").Value2 = 1.87
$ws.Range("V54").Value2 = 1.87
$ws.Range("M55").Value2 = 1.04
$ws.Range("O55").Value2 = 1.2
$ws.Range("U55").Value2 = 1.58
$ws.Range("G56").Value2 = 1.5
$ws.Range("H56").Value2 = 4.1
$ws.Range("I56").Value2 = 6.5
$ws.Range("J56").Value2 = 2
$ws.Range("K56").Value2 = 2.4
$ws.Range("M56").Value2 = 1.04
$ws.Range("O56").Value2 = 1.22
$ws.Range("U56").Value2 = 1.77
$ws.Range("V56").Value2 = 1.92
$ws.Range("Z56").Value2 = 10
$ws.Range("AD56").Value2 = 8
$ws.Range("AM56").Value2 = 51
$ws.Range("AX56").Value2 = 34
$ws.Range("AZ56").Value2 = 126
$ws.Range("M57").Value2 = 1.02
$ws.Range("O57").Value2 = 1.13
$ws.Range("G61").Value2 = 2.05
$ws.Range("U61").Value2 = 1.37
$ws.Range("G62").Value2 = 1.57
$ws.Range("N62").Value2 = 26
$ws.Range("U62").Value2 = 1.33
$ws.Range("G63").Value2 = 2.1
$ws.Range("U63").Value2 = 1.54
$ws.Range("G64").Value2 = 1.8
$ws.Range("U64").Value2 = 1.54
$ws.Range("G65").Value2 = 2.63
$ws.Range("I65").Value2 = 2.55
$ws.Range("J65").Value2 = 3.2
$ws.Range("L65").Value2 = 3.1
$ws.Range("U65").Value2 = 1.54
$ws.Range("Z65").Value2 = 26
$ws.Range("AA65").Value2 = 19
$ws.Range("AK65").Value2 = 26
$ws.Range("AN65").Value2 = 4.75
$ws.Range("Q69").Value2 = 1.67
$ws.Range("U69").Value2 = 1.62
$ws.Range("Q70").Value2 = 1.57
$ws.Range("U70").Value2 = 1.5
$ws.Range("R71").Value2 = 1.65
$ws.Range("U71").Value2 = 1.83
$ws.Range("V71").Value2 = 1.83
$ws.Range("R72").Value2 = 1.67
$ws.Range("R73").Value2 = 1.67
$ws.Range("Q74").Value2 = 1.73
$ws.Range("R74").Value2 = 2.08
$ws.Range("R75").Value2 = 1.48
$ws.Range("Q76").Value2 = 1.8
$ws.Range("R77").Value2 = 1.7
$ws.Range("M80").Value2 = 1.08
$ws.Range("O80").Value2 = 1.4
$ws.Range("U80").Value2 = 1.83
$ws.Range("V80").Value2 = 1.83

# --- Insert a new row at position 85 (Serbia - Super Liga: Cukaricki vs Napredak) ---
# This shifts the existing rows 85-91 down to 86-92.
$ws.Rows(85).Insert()

$r85 = New-Object 'object[,]' 1,56
$r85[0,0] = '6329OsCa'
$r85[0,1] = '18/10/2024'
$r85[0,2] = '13:00'
$r85[0,3] = 'SERBIA - SUPER LIGA'
$r85[0,4] = 'Cukaricki'
$r85[0,5] = 'Napredak'
$r85[0,6] = 1.42
$r85[0,7] = 4.5
$r85[0,8] = 6.5
$r85[0,9] = 1.88
$r85[0,10] = 2.42
$r85[0,11] = 6
$r85[0,12] = 1.04
$r85[0,13] = 8.75
$r85[0,14] = 1.2
$r85[0,15] = 4.05
$r85[0,16] = 1.62
$r85[0,17] = 2.22
$r85[0,18] = 1.31
$r85[0,19] = 3.15
$r85[0,20] = 1.8
$r85[0,21] = 1.9
$r85[0,22] = 7.9
$r85[0,23] = 7.3
$r85[0,24] = 8.25
$r85[0,25] = 9.75
$r85[0,26] = 11
$r85[0,27] = 24
$r85[0,28] = 8.75
$r85[0,29] = 9
$r85[0,30] = 18
$r85[0,31] = 75
$r85[0,32] = 500
$r85[0,33] = 19
$r85[0,34] = 40
$r85[0,35] = 20
$r85[0,36] = 150
$r85[0,37] = 65
$r85[0,38] = 60
$r85[0,39] = 3.35
$r85[0,40] = 6.4
$r85[0,41] = 15.5
$r85[0,42] = 18
$r85[0,43] = 40
$r85[0,44] = 200
$r85[0,45] = 3.15
$r85[0,46] = 8
$r85[0,47] = 70
$r85[0,48] = 7.9
$r85[0,49] = 37
$r85[0,50] = 37
$r85[0,51] = 250
$r85[0,52] = 250
$r85[0,53] = 450
$r85[0,54] = 51
$r85[0,55] = 51
$ws.Range("A85:BD85").Value = $r85

# --- Append a new row 93 (Ukraine - Premier League: Kolos Kovalivka vs Shakhtar Donetsk) ---
$r93 = New-Object 'object[,]' 1,56
$r93[0,0] = 'WdfySwEF'
$r93[0,1] = '18/10/2024'
$r93[0,2] = '12:00'
$r93[0,3] = 'UKRAINE - PREMIER LEAGUE'
$r93[0,4] = 'Kolos Kovalivka'
$r93[0,5] = 'Shakhtar Donetsk'
$r93[0,6] = 9.75
$r93[0,7] = 4.75
$r93[0,8] = 1.29
$r93[0,9] = 8.25
$r93[0,10] = 2.35
$r93[0,11] = 1.75
$r93[0,12] = 1.01
$r93[0,13] = 12.5
$r93[0,14] = 1.25
$r93[0,15] = 3.25
$r93[0,16] = 1.78
$r93[0,17] = 1.93
$r93[0,18] = 1.35
$r93[0,19] = 3.02
$r93[0,20] = 2.25
$r93[0,21] = 1.5
$r93[0,22] = 21
$r93[0,23] = 70
$r93[0,24] = 32
$r93[0,25] = 300
$r93[0,26] = 150
$r93[0,27] = 150
$r93[0,28] = 10.5
$r93[0,29] = 10
$r93[0,30] = 29
$r93[0,31] = 175
$r93[0,32] = 67
$r93[0,33] = 5.9
$r93[0,34] = 5.5
$r93[0,35] = 9
$r93[0,36] = 7.2
$r93[0,37] = 12
$r93[0,38] = 37
$r93[0,39] = 9.75
$r93[0,40] = 65
$r93[0,41] = 65
$r93[0,42] = 500
$r93[0,43] = 500
$r93[0,44] = 67
$r93[0,45] = 2.65
$r93[0,46] = 10
$r93[0,47] = 120
$r93[0,48] = 2.92
$r93[0,49] = 5.5
$r93[0,50] = 18.5
$r93[0,51] = 15
$r93[0,52] = 50
$r93[0,53] = 300
$r93[0,54] = 81
$r93[0,55] = 81
$ws.Range("A93:BD93").Value = $r93
